# Update 20 July 2024
# Applies the payroll-ledger updates described by the commit:
#  - New "Keterangan" (notes) status text on a batch of existing rows
#  - A couple of date/amount corrections on a few existing rows
#  - Eight brand-new payment rows appended at the bottom (174-181)
#  - Table1 / autofilter / used-range resized to match
#  - Selection moved to reflect where the user was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New "Keterangan" notes added to rows that didn't have one yet
# ---------------------------------------------------------------------------
$ws.Range("N116").Value = "gabung tgl 01 juni 2024"
$ws.Range("N133").Value = "salah total, ditagih tgl 23 juni 2024"
$ws.Range("N140").Value = "gabunt tgl 2 juli 2024"
$ws.Range("N156").Value = "termasuk tgl 1 juli"
$ws.Range("N157").Value = "tesmasuk tgl 1 juli"
$ws.Range("N159").Value = "gabung tgl 16 juli 2024"
$ws.Range("N160").Value = "belum lunas"
$ws.Range("N168").Value = "gabung 21 juli 2024"

# ---------------------------------------------------------------------------
# 2. "Keterangan" notes updated from "belum lunas" (still unpaid) to the
#    actual status now that it is resolved
# ---------------------------------------------------------------------------
$ws.Range("N109").Value = "gabung tgl 03 juli 2024"
$ws.Range("N118").Value = "gabung tgl 07 juni 2024"
$ws.Range("N119").Value = "gabung tgl 25 juni 2024"
$ws.Range("N129").Value = "gabung tgl 16 juni 2024"
$ws.Range("N141").Value = "gabung tgl 1 juli 2024"
$ws.Range("N144").Value = "gabung tgl 06 juli 2024"
$ws.Range("N169").Value = "gabung tgl 13 juli 2024"
$ws.Range("N170").Value = "gabung tgl 21 juli 2024"
$ws.Range("N173").Value = "gabung tgl 13 juli 2024"

# ---------------------------------------------------------------------------
# 3. Small corrections on a few existing rows
# ---------------------------------------------------------------------------
$ws.Range("D149").Value = 45451
$ws.Range("E151").Value = 45472

$ws.Range("E160").Value = 45465
$ws.Range("F160").Value = 15148000
$ws.Range("G160").Value = 15148000

$ws.Range("M173").Value = 0

# ---------------------------------------------------------------------------
# 4. Eight brand-new payment rows appended at the bottom of the table
# ---------------------------------------------------------------------------

# Row 174
$ws.Range("A174").Value = 45486
$ws.Range("B174").Value = "Subadi"
$ws.Range("D174").Value = 45467
$ws.Range("E174").Value = 45467
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 26200000
$ws.Range("I174").Formula = "=G174-F174"
$ws.Range("J174").Value = 10000000
$ws.Range("K174").Formula = "=G174-J174"
$ws.Range("L174").Formula = "=G174-J174+H174"
$ws.Range("M174").Value = 0
$ws.Range("N174").Value = "gabung tgl 15 juli 2024"

# Row 175
$ws.Range("B175").Value = "Amin"
$ws.Range("D175").Value = 45443
$ws.Range("E175").Value = 45456
$ws.Range("F175").Value = 2675000
$ws.Range("G175").Value = 2675000
$ws.Range("H175").Value = 565000
$ws.Range("I175").Formula = "=G175-F175"
$ws.Range("J175").Value = 2675000
$ws.Range("K175").Formula = "=G175-J175"
$ws.Range("L175").Formula = "=G175-J175+H175"
$ws.Range("M175").Value = 100000
$ws.Range("N175").Value = "belum lunas"

# Row 176
$ws.Range("A176").Value = 45488
$ws.Range("B176").Value = "Subadi"
$ws.Range("D176").Value = 45468
$ws.Range("E176").Value = 45485
$ws.Range("F176").Value = 2992000
$ws.Range("G176").Value = 2992000
$ws.Range("H176").Value = 16200000
$ws.Range("I176").Formula = "=G176-F176"
$ws.Range("J176").Value = 2565000
$ws.Range("K176").Formula = "=G176-J176"
$ws.Range("L176").Formula = "=G176-J176+H176"
$ws.Range("M176").Value = 45000
$ws.Range("N176").Value = "belum lunas"

# Row 177
$ws.Range("B177").Value = "Fatoni"
$ws.Range("D177").Value = 45473
$ws.Range("E177").Value = 45486
$ws.Range("F177").Value = 6732000
$ws.Range("G177").Value = 6732000
$ws.Range("H177").Value = 0
$ws.Range("I177").Formula = "=G177-F177"
$ws.Range("J177").Value = 6031000
$ws.Range("K177").Formula = "=G177-J177"
$ws.Range("L177").Formula = "=G177-J177+H177"
$ws.Range("M177").Value = 100000

# Row 178
$ws.Range("A178").Value = 45489
$ws.Range("B178").Value = "Sugeng"
$ws.Range("D178").Value = 45474
$ws.Range("E178").Value = 45487
$ws.Range("F178").Value = 5359000
$ws.Range("G178").Value = 5359000
$ws.Range("H178").Value = 1733000
$ws.Range("I178").Formula = "=G178-F178"
$ws.Range("J178").Value = 6500000
$ws.Range("K178").Formula = "=G178-J178"
$ws.Range("L178").Formula = "=G178-J178+H178"
$ws.Range("M178").Value = 100000
$ws.Range("N178").Value = "belum lunas"

# Row 179
$ws.Range("B179").Value = "Manggi"
$ws.Range("D179").Value = 45466
$ws.Range("E179").Value = 45478
$ws.Range("F179").Value = 13163000
$ws.Range("G179").Value = 13163000
$ws.Range("H179").Value = 1754000
$ws.Range("I179").Formula = "=G179-F179"
$ws.Range("J179").Value = 14917000
$ws.Range("K179").Formula = "=G179-J179"
$ws.Range("L179").Formula = "=G179-J179+H179"
$ws.Range("M179").Value = 100000
$ws.Range("N179").Value = "lunas"

# Row 180
$ws.Range("A180").Value = 45494
$ws.Range("B180").Value = "Dede"
$ws.Range("D180").Value = 45480
$ws.Range("E180").Value = 45492
$ws.Range("F180").Value = 12105000
$ws.Range("G180").Value = 12105000
$ws.Range("H180").Value = 68000
$ws.Range("I180").Formula = "=G180-F180"
$ws.Range("J180").Value = 12105000
$ws.Range("K180").Formula = "=G180-J180"
$ws.Range("L180").Formula = "=G180-J180+H180"
$ws.Range("M180").Value = 100000
$ws.Range("N180").Value = "belum lunas"

# Row 181
$ws.Range("B181").Value = "Andre"
$ws.Range("D181").Value = 45478
$ws.Range("E181").Value = 45492
$ws.Range("F181").Value = 2916000
$ws.Range("G181").Value = 2916000
$ws.Range("H181").Value = 220000
$ws.Range("I181").Formula = "=G181-F181"
$ws.Range("J181").Value = 2916000
$ws.Range("K181").Formula = "=G181-J181"
$ws.Range("L181").Formula = "=G181-J181+H181"
$ws.Range("M181").Value = 100000
$ws.Range("N181").Value = "belum lunas"

# ---------------------------------------------------------------------------
# 5. Grow Table1 (and its AutoFilter) to cover the newly-added rows
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:N181"))

# ---------------------------------------------------------------------------
# 6. Move the view / selection to where the user last left off
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 141
$ws.Range("E152").Select()
